$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-36 (new order), pipe-delimited:
# A | B(Datum) | C(Forandrad) | G(Area) | H(Fridlysta) | I(Signalarter) | J(NT) | K(VU) | L(EN) | M(CR) | N(RE) | O(Rodlistade) | P(Hotade) | Q(Alla arter) | R(Artnamn)
$csvData = @"
A 24087-2025|45795|46078|30.3|0|0|1|0|0|0|0|1|0|1|Rödlånke
A 24088-2025|45795|46078|2.8|1|0|0|0|0|0|0|0|0|1|Blåsippa
A 24036-2025|45795|46078|37.1|1|0|0|0|0|0|0|0|0|1|Kopparödla
A 2980-2026|46038.62965277778|46078|1.1|1|0|0|0|0|0|0|0|0|1|Blåsippa
A 18791-2021|44307.58013888889|46078|1|0|0|0|0|0|0|0|0|0|0|
A 54967-2022|44886.4831712963|46078|0.5|0|0|0|0|0|0|0|0|0|0|
A 43053-2021|44431|46078|3|0|0|0|0|0|0|0|0|0|0|
A 43056-2021|44431|46078|1|0|0|0|0|0|0|0|0|0|0|
A 21487-2021|44320|46078|1|0|0|0|0|0|0|0|0|0|0|
A 73020-2021|44550|46078|2.7|0|0|0|0|0|0|0|0|0|0|
A 59231-2024|45637.58472222222|46078|1.3|0|0|0|0|0|0|0|0|0|0|
A 28416-2024|45477.62978009259|46078|0.2|0|0|0|0|0|0|0|0|0|0|
A 21421-2021|44316|46078|0.6|0|0|0|0|0|0|0|0|0|0|
A 17491-2024|45415.50266203703|46078|6.2|0|0|0|0|0|0|0|0|0|0|
A 50864-2022|44867.56143518518|46078|3.3|0|0|0|0|0|0|0|0|0|0|
A 23503-2025|45795|46078|14.1|0|0|0|0|0|0|0|0|0|0|
A 270-2025|45660.48087962963|46078|8.9|0|0|0|0|0|0|0|0|0|0|
A 49633-2024|45596.59559027778|46078|0.8|0|0|0|0|0|0|0|0|0|0|
A 24086-2025|45795|46078|0.7|0|0|0|0|0|0|0|0|0|0|
A 4422-2024|45327.45375|46078|4.5|0|0|0|0|0|0|0|0|0|0|
A 24212-2023|45076|46078|5.8|0|0|0|0|0|0|0|0|0|0|
A 51434-2025|45949|46078|2.8|0|0|0|0|0|0|0|0|0|0|
A 50239-2022|44865|46078|13.2|0|0|0|0|0|0|0|0|0|0|
A 46579-2024|45582.75018518518|46078|3|0|0|0|0|0|0|0|0|0|0|
A 52965-2025|45956|46078|0.6|0|0|0|0|0|0|0|0|0|0|
A 52960-2025|45956|46078|1.7|0|0|0|0|0|0|0|0|0|0|
A 52888-2025|45957.56943287037|46078|2|0|0|0|0|0|0|0|0|0|0|
A 28409-2024|45477.62280092593|46078|0.4|0|0|0|0|0|0|0|0|0|0|
A 57001-2025|45977|46078|2|0|0|0|0|0|0|0|0|0|0|
A 4780-2022|44592.62657407407|46078|0.9|0|0|0|0|0|0|0|0|0|0|
A 17492-2024|45415.50709490741|46078|7.7|0|0|0|0|0|0|0|0|0|0|
A 46587-2024|45582.76763888889|46078|2.8|0|0|0|0|0|0|0|0|0|0|
A 46588-2024|45582.77137731481|46078|1|0|0|0|0|0|0|0|0|0|0|
A 28418-2024|45477.62978009259|46078|0.2|0|0|0|0|0|0|0|0|0|0|
A 49634-2024|45596.59591435185|46078|0.7|0|0|0|0|0|0|0|0|0|0|
"@

$baseUrl = "https://klasma.github.io/Logging_1427"

$lines = $csvData -split "`r`n|`n"
$rowNum = 2
foreach ($line in $lines) {
    $f = $line -split "\|"
    $beteckning = $f[0]
    $datum = [double]$f[1]
    $forandrad = [double]$f[2]
    $area = [double]$f[3]
    $fridlysta = [double]$f[4]
    $signalarter = [double]$f[5]
    $nt = [double]$f[6]
    $vu = [double]$f[7]
    $en = [double]$f[8]
    $cr = [double]$f[9]
    $re = [double]$f[10]
    $rodlistade = [double]$f[11]
    $hotade = [double]$f[12]
    $allaArter = [double]$f[13]
    $artnamn = $f[14]

    $ws.Cells.Item($rowNum, 1).Value2 = $beteckning
    $ws.Cells.Item($rowNum, 2).Value2 = $datum
    $ws.Cells.Item($rowNum, 3).Value2 = $forandrad
    $ws.Cells.Item($rowNum, 7).Value2 = $area
    $ws.Cells.Item($rowNum, 8).Value2 = $fridlysta
    $ws.Cells.Item($rowNum, 9).Value2 = $signalarter
    $ws.Cells.Item($rowNum, 10).Value2 = $nt
    $ws.Cells.Item($rowNum, 11).Value2 = $vu
    $ws.Cells.Item($rowNum, 12).Value2 = $en
    $ws.Cells.Item($rowNum, 13).Value2 = $cr
    $ws.Cells.Item($rowNum, 14).Value2 = $re
    $ws.Cells.Item($rowNum, 15).Value2 = $rodlistade
    $ws.Cells.Item($rowNum, 16).Value2 = $hotade
    $ws.Cells.Item($rowNum, 17).Value2 = $allaArter

    if ($artnamn -ne "") {
        $ws.Cells.Item($rowNum, 18).Value2 = $artnamn

        $ws.Cells.Item($rowNum, 19).Formula = '=HYPERLINK("' + $baseUrl + '/artfynd/' + $beteckning + ' artfynd.xlsx", "' + $beteckning + '")'
        $ws.Cells.Item($rowNum, 20).Formula = '=HYPERLINK("' + $baseUrl + '/kartor/' + $beteckning + ' karta.png", "' + $beteckning + '")'
        $ws.Cells.Item($rowNum, 22).Formula = '=HYPERLINK("' + $baseUrl + '/klagomål/' + $beteckning + ' FSC-klagomål.docx", "' + $beteckning + '")'
        $ws.Cells.Item($rowNum, 23).Formula = '=HYPERLINK("' + $baseUrl + '/klagomålsmail/' + $beteckning + ' FSC-klagomål mail.docx", "' + $beteckning + '")'
        $ws.Cells.Item($rowNum, 24).Formula = '=HYPERLINK("' + $baseUrl + '/tillsyn/' + $beteckning + ' tillsynsbegäran.docx", "' + $beteckning + '")'
        $ws.Cells.Item($rowNum, 25).Formula = '=HYPERLINK("' + $baseUrl + '/tillsynsmail/' + $beteckning + ' tillsynsbegäran mail.docx", "' + $beteckning + '")'
    }

    $rowNum++
}
